$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.084.87'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '1.650.02'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  -0.50%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5197'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.82%  '
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  -1.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06278'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07807'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.456'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.30%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.635.29'
$ws.Range('E13').Value = '  -1.95%  '
$ws.Range('D14').Value = '1.876.88'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5542'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '0.0₅7984'
$ws.Range('E16').Value = '  -2.81%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '26.079.17'
$ws.Range('E18').Value = '  -0.96%  '
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.625'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '194.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.937'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.007'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1204'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.178'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.86'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.57%  '
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05598'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.263'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.471'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.373'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.590'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.798'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9476'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.404'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5645'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.89%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.972'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01578'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('D41').Value = '1.060.79'
$ws.Range('E41').Value = '  +0.87%  '
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8386'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.26%  '
$ws.Range('D45').Value = '1.788.56'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05351'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.61%  '
$ws.Range('E49').Value = '  -0.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4334'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.936'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.27%  '
